$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data updates (PERIOD TO EXPIRE / LAST UPDATE) ---------------------
# H3 / H4: period-to-expire numbers move from -40 to -48.
$ws.Range("H3").Value = -48
$ws.Range("H4").Value = -48

# I3 / I4: last-update date strings move from 08-Sep-2025 to 16-Sep-2025.
# These are stored as literal text (not real dates) in the sheet, so force
# the cell to text first, then restore the original number format (which
# lives on H3/H4, same row/style) via a formats-only paste so the cell's
# style index is unaffected by the temporary text format.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
$ws.Range("H3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "16-Sep-2025"
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Styling: header rows get white bold text (instead of plain bold) --
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215
